$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Tear down the old layout: remove merges + clear all cell content/styles
#    so the sheet can be rebuilt cleanly in the new shape.
# ---------------------------------------------------------------------------
$ws.Cells.UnMerge()
$ws.Cells.Clear()

# ---------------------------------------------------------------------------
# 2. Column setup: add the new column J (used by the "Следующий Месяц" block)
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 8

# ---------------------------------------------------------------------------
# 3. Header row 1 - one label per 2-column block (the pair is merged below)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Сегодня "
$ws.Range("C1").Value = "Завтра "
$ws.Range("E1").Value = "Неделя"
$ws.Range("G1").Value = "Месяц"
$ws.Range("I1").Value = "Следующий Месяц"

# ---------------------------------------------------------------------------
# 4. Header row 2 - "Задача" / "Ответственный" pairs repeated under each block
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Задача"
$ws.Range("B2").Value = "Ответственный "
$ws.Range("C2").Value = "Задача"
$ws.Range("D2").Value = "Ответственный "
$ws.Range("E2").Value = "Задача"
$ws.Range("F2").Value = "Ответственный "
$ws.Range("G2").Value = "Задача"
$ws.Range("H2").Value = "Ответственный "

# ---------------------------------------------------------------------------
# 5. Data rows
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Разработать шаблон excel файла"
$ws.Range("B3").Value = "Копытов П.Е."
$ws.Range("C3").Value = "Разработать шаблон нового excel файла "
$ws.Range("D3").Value = "Копытов П.Е."
$ws.Range("E3").Value = "Создать парсер "
$ws.Range("F3").Value = "Копытоа П.Е."

$ws.Range("A4").Value = "описать структуру работы и программы"
$ws.Range("B4").Value = "Копытов П.Е."
$ws.Range("C4").Value = "Доработать описание"
$ws.Range("D4").Value = "Копытов П.Е."
$ws.Range("E4").Value = "Выбрать наиболее подходящий шаблон"
$ws.Range("F4").Value = "Копытоа П.Е."

$ws.Range("C5").Value = "фосфор"
$ws.Range("D5").Value = "Копт"
$ws.Range("E5").Value = "Тестирование парсера"
$ws.Range("F5").Value = "Копытоа П.Е."

$ws.Range("C6").Value = "ФПЙ"
$ws.Range("D6").Value = "циц"
$ws.Range("E6").Value = "УКТФИ"
$ws.Range("F6").Value = "Копытоа П.Е."

$ws.Range("C7").Value = "ицыиы"
$ws.Range("D7").Value = "ыиыу"
$ws.Range("E7").Value = "УКИВ"
$ws.Range("F7").Value = "Кымиы"

$ws.Range("C8").Value = 4
$ws.Range("D8").Value = "циы"
$ws.Range("E8").Value = "ПРЕИУ"
$ws.Range("F8").Value = "циц"

$ws.Range("C9").Value = 34

# ---------------------------------------------------------------------------
# 6. Styling: row 1 centred, row 2 bold (reuses the workbook's existing
#    style slots, same as the original template)
# ---------------------------------------------------------------------------
$ws.Range("A1:J1").HorizontalAlignment = -4108
$ws.Range("A2:H2").Font.Bold = $true

# ---------------------------------------------------------------------------
# 7. Re-merge the header blocks (shifted one column left vs. the old layout,
#    plus the new trailing I1:J1 block)
# ---------------------------------------------------------------------------
$ws.Range("A1:B1").Merge()
$ws.Range("C1:D1").Merge()
$ws.Range("E1:F1").Merge()
$ws.Range("G1:H1").Merge()
$ws.Range("I1:J1").Merge()

# ---------------------------------------------------------------------------
# 8. Selection matches the saved cursor position in the edited workbook
# ---------------------------------------------------------------------------
[void]$ws.Range("A2").Select()
